# Client docx v2 with query
#
# Appends a centered "Query (operazioni 1,2,3,4)" title and four
# bold-numbered query descriptions at the end of the document, followed
# by a trailing empty paragraph (mirroring the existing trailing
# <w:tab/> paragraph that precedes them).

$d = $word.ActiveDocument

# Insert six brand-new, still-empty paragraphs right away (before any
# formatting is applied to any of them). Because each new paragraph
# inherits the mark formatting of the paragraph that exists at the
# moment it is created, doing all the insertions first -- while the
# preceding paragraph is still the plain "<w:tab/>" one -- keeps every
# new paragraph free of inherited bold/center/size formatting. We then
# go back and format only the specific paragraphs/ranges that need it.
$count = $d.Paragraphs.Count
for ($i = 0; $i -lt 6; $i++) {
    $d.Paragraphs.Last.Range.InsertParagraphAfter()
}

$titleIndex = $count + 1
$p1Index    = $count + 2
$p2Index    = $count + 3
$p3Index    = $count + 4
$p4Index    = $count + 5
# $count + 6 is the final trailing empty paragraph -- left untouched.

# --- "Query (operazioni 1,2,3,4)" title -----------------------------------
$title = $d.Paragraphs.Item($titleIndex)
$title.Range.Text = "Query (operazioni 1,2,3,4)"
$title.Alignment = 1
$title.Range.Font.Bold = 1
$title.Range.Font.Size = 18
$title.Range.Font.SizeBi = 18

# --- helper to build "<n> = <text>" paragraphs with a bold leading number -
function Set-QueryParagraph($paragraphIndex, [string]$number, [string]$body) {
    $p = $d.Paragraphs.Item($paragraphIndex)
    $p.Range.Text = $number + $body
    $start = $p.Range.Start
    $boldRange = $d.Range($start, $start + $number.Length)
    $boldRange.Font.Bold = 1
}

$quote = [char]0x201C
$unquote = [char]0x201D

Set-QueryParagraph $p1Index "1" (" = stampare quantità dei libri " + $quote + "Ultimi Arrivi" + $unquote + " nel reparto " + $quote + "fumetti" + $unquote)

Set-QueryParagraph $p2Index "2" " = elenco dei libri scontati presenti in tutti i reparti in ordine crescente per sconto (da quelli meno a quelli più scontati)"

Set-QueryParagraph $p3Index "3" " = elenco libri archiviati all'interno di un periodo definito da due date inserite in input"

Set-QueryParagraph $p4Index "4" " = dato il codice di un carrello elenco dei titoli dei libri acquistati con il rispettivo numero copie e username dell'utente associato a quel carrell"
